$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("A1").Value = "neighbourhood"

# Sort the data (A2:B129) by averagePrice (column B) ascending
$rangeToSort = $ws.Range("A1:B129")
$sortField = $ws.Range("B1:B129")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortField, 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($rangeToSort)
$ws.Sort.Header = 1
$ws.Sort.Apply()
